# recursive-data.pptx edit
#
# Per the commit diff, slide 2's title ("Recursive " + "Data" split across two
# runs, with a trailing empty <a:endParaRPr>) is re-typed as a single run
# "Recursive Data" with no trailing endParaRPr.
#
# (The diff's other hunks only touch the legacy VML `spid="_x0000_sNNNN"`
# fallback ids on the embedded Equation OLE objects on slides 10 & 14 -
# an internal counter PowerPoint itself assigns when it resaves the file
# and which isn't exposed anywhere in the PowerPoint object model - plus a
# purely cosmetic namespace-declaration reshuffle of the slide 14 transition's
# mc:AlternateContent/mc:Fallback wrapper. Neither is reachable from
# COM/VBA automation, and this deck's transition markup already matches the
# target layout, so there is nothing further to change there.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)

# Clear the existing run(s) first so the engine rebuilds the paragraph from
# scratch as a single run, rather than just patching the text of the existing
# runs (which would leave the "Recursive "/"Data" split and the stray
# endParaRPr in place).
$sh.TextFrame.TextRange.Delete()
$sh.TextFrame.TextRange.Text = "Recursive Data"
